$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Activate()
$ws.Range("E8").Value = "GIT UPDATE"
$ws.Range("E8").Select()
